$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2038834951456311
$ws.Range("C2").Value = 0.540453074433657
$ws.Range("J2").Value = 0.009708737864077669
$ws.Range("O2").Value = 0.003236245954692557
$ws.Range("P2").Value = 0.1682847896440129
$ws.Range("S2").Value = 0.07443365695792881
$ws.Range("B3").Value = 0.02352941176470588
$ws.Range("J3").Value = 0.05294117647058823
$ws.Range("P3").Value = 0.7705882352941177
$ws.Range("S3").Value = 0.1529411764705882
$ws.Range("J4").Value = 0.05555555555555555
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2777777777777778
$ws.Range("B6").Value = 0.04661016949152542
$ws.Range("D6").Value = 0.02542372881355932
$ws.Range("F6").Value = 0.05084745762711865
$ws.Range("J6").Value = 0.2457627118644068
$ws.Range("O6").Value = 0.0211864406779661
$ws.Range("Q6").Value = 0.1567796610169492
$ws.Range("R6").Value = 0.07203389830508475
$ws.Range("S6").Value = 0.3813559322033898
$ws.Range("B7").Value = 0.1625
$ws.Range("D7").Value = 0.01875
$ws.Range("F7").Value = 0.05625
$ws.Range("J7").Value = 0.175
$ws.Range("O7").Value = 0.0125
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.075
$ws.Range("S7").Value = 0.35
$ws.Range("B8").Value = 0.08884297520661157
$ws.Range("D8").Value = 0.02066115702479339
$ws.Range("F8").Value = 0.07231404958677685
$ws.Range("J8").Value = 0.134297520661157
$ws.Range("O8").Value = 0.03099173553719008
$ws.Range("Q8").Value = 0.1900826446280992
$ws.Range("R8").Value = 0.08884297520661157
$ws.Range("S8").Value = 0.3739669421487603
$ws.Range("B9").Value = 0.09844559585492228
$ws.Range("D9").Value = 0.03626943005181347
$ws.Range("E9").Value = 0.005181347150259068
$ws.Range("F9").Value = 0.05699481865284974
$ws.Range("J9").Value = 0.1191709844559585
$ws.Range("O9").Value = 0.01036269430051814
$ws.Range("Q9").Value = 0.2124352331606218
$ws.Range("R9").Value = 0.08290155440414508
$ws.Range("S9").Value = 0.3782383419689119
$ws.Range("B10").Value = 0.100211714890614
$ws.Range("D10").Value = 0.02117148906139732
$ws.Range("F10").Value = 0.068454481298518
$ws.Range("J10").Value = 0.1362032462949894
$ws.Range("O10").Value = 0.01905434015525759
$ws.Range("Q10").Value = 0.2491178546224418
$ws.Range("R10").Value = 0.07339449541284404
$ws.Range("S10").Value = 0.3323923782639379
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.0946969696969697
$ws.Range("K11").Value = 0.1818181818181818
$ws.Range("L11").Value = 0.5909090909090909
$ws.Range("S11").Value = 0.007575757575757576
$ws.Range("G12").Value = 0.7018633540372671
$ws.Range("J12").Value = 0.2360248447204969
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("S12").Value = 0.01863354037267081
$ws.Range("G13").Value = 0.5151515151515151
$ws.Range("J13").Value = 0.4545454545454545
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("F15").Value = 0.024
$ws.Range("H15").Value = 0.188
$ws.Range("I15").Value = 0.068
$ws.Range("J15").Value = 0.34
$ws.Range("K15").Value = 0.06
$ws.Range("M15").Value = 0.016
$ws.Range("O15").Value = 0.06
$ws.Range("S15").Value = 0.244
$ws.Range("F16").Value = 0.03846153846153846
$ws.Range("H16").Value = 0.1730769230769231
$ws.Range("I16").Value = 0.09134615384615384
$ws.Range("J16").Value = 0.4086538461538461
$ws.Range("K16").Value = 0.07692307692307693
$ws.Range("M16").Value = 0.004807692307692308
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.0576923076923077
$ws.Range("S16").Value = 0.1442307692307692
$ws.Range("F17").Value = 0.02022058823529412
$ws.Range("H17").Value = 0.2022058823529412
$ws.Range("I17").Value = 0.08823529411764706
$ws.Range("J17").Value = 0.4172794117647059
$ws.Range("K17").Value = 0.08088235294117647
$ws.Range("M17").Value = 0.01470588235294118
$ws.Range("N17").Value = 0.001838235294117647
$ws.Range("O17").Value = 0.07720588235294118
$ws.Range("S17").Value = 0.09742647058823529
$ws.Range("F18").Value = 0.01036269430051814
$ws.Range("H18").Value = 0.1606217616580311
$ws.Range("I18").Value = 0.09326424870466321
$ws.Range("J18").Value = 0.4818652849740933
$ws.Range("K18").Value = 0.09326424870466321
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.07253886010362694
$ws.Range("S18").Value = 0.07772020725388601
$ws.Range("F19").Value = 0.01402640264026403
$ws.Range("H19").Value = 0.2194719471947195
$ws.Range("I19").Value = 0.07425742574257425
$ws.Range("J19").Value = 0.4018151815181518
$ws.Range("K19").Value = 0.09818481848184818
$ws.Range("M19").Value = 0.01485148514851485
$ws.Range("N19").Value = 0.0008250825082508251
$ws.Range("O19").Value = 0.07178217821782178
$ws.Range("S19").Value = 0.1047854785478548
